$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '98.674.47'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.30%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.317.11'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '256.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '623.54'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.92%  '
$ws.Range('E7').Value = '  +29.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.404'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.70%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.912'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +16.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.317.40'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '39.17'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +11.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '98.449.66'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000249'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.938.44'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.49'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.316.01'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.49'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.29'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +8.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '483.73'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.44'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000205'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.62'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '88.72'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.486.65'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.81%  '
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.293'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +23.23%  '
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.189'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.15%  '
$ws.Range('E32').Value = '  +8.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '10.23'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +11.32%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '27.99'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.27'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.149'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.95'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.42%  '
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.463'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.06%  '
$ws.Range('B40').Value = 'WhiteBITCoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '24.83'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '492.50'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.24'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.05%  '
$ws.Range('B43').Value = 'MantraDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.63'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.796'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.25%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '158.32'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.39'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +16.81%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.92'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.71%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.846'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.65'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.70%  '
